$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2026-02-28 -> 2026-03-01, serial 46081 -> 46082) for every data row (2-37).
for ($row = 2; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
